$d = $word.ActiveDocument

# Update the date heading (first paragraph).
$d.Paragraphs.Item(1).Range.Text = "2023-10-24 Tuesday"

# Update the 25 answer cells in the 5x5 practice table.
# The table has 20 rows total (data rows at 1,5,9,13,17; the rows
# between them are blank spacer rows), 5 columns.
$t = $d.Tables.Item(1)

$newValues = @(
    @(1,  1, "40÷4=10, 0"),
    @(1,  2, "36÷2=18, 0"),
    @(1,  3, "52÷7=7, 3"),
    @(1,  4, "96÷9=10, 6"),
    @(1,  5, "37÷4=9, 1"),

    @(5,  1, "21÷2=10, 1"),
    @(5,  2, "61÷3=20, 1"),
    @(5,  3, "80÷4=20, 0"),
    @(5,  4, "94÷4=23, 2"),
    @(5,  5, "69÷8=8, 5"),

    @(9,  1, "97÷7=13, 6"),
    @(9,  2, "18÷9=2, 0"),
    @(9,  3, "54÷2=27, 0"),
    @(9,  4, "25÷8=3, 1"),
    @(9,  5, "71÷2=35, 1"),

    @(13, 1, "91÷3=30, 1"),
    @(13, 2, "53÷5=10, 3"),
    @(13, 3, "65÷5=13, 0"),
    @(13, 4, "37÷9=4, 1"),
    @(13, 5, "19÷3=6, 1"),

    @(17, 1, "51÷7=7, 2"),
    @(17, 2, "27÷5=5, 2"),
    @(17, 3, "49÷2=24, 1"),
    @(17, 4, "10÷4=2, 2"),
    @(17, 5, "90÷3=30, 0")
)

foreach ($entry in $newValues) {
    $rowIdx = $entry[0]
    $colIdx = $entry[1]
    $newText = $entry[2]
    $cell = $t.Cell($rowIdx, $colIdx)
    $cell.Range.Text = $newText
}
